$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s_D2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.050.28'
$ws.Range("D2").Style = $s_D2
$s_E2 = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("E2").Style = $s_E2
$s_D3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.172.21'
$ws.Range("D3").Style = $s_D3
$s_E3 = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("E3").Style = $s_E3
$s_D4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = $s_D4
$s_E4 = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E4").Style = $s_E4
$s_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.81'
$ws.Range("D5").Style = $s_D5
$s_E5 = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("E5").Style = $s_E5
$s_E6 = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("E6").Style = $s_E6
$s_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '66.32'
$ws.Range("D7").Style = $s_D7
$s_E7 = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -5.44%  '
$ws.Range("E7").Style = $s_E7
$s_E8 = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E8").Style = $s_E8
$s_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.559'
$ws.Range("D9").Style = $s_D9
$s_E9 = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("E9").Style = $s_E9
$s_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '60.04'
$ws.Range("D10").Style = $s_D10
$s_E10 = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.49%  '
$ws.Range("E10").Style = $s_E10
$s_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0925'
$ws.Range("D11").Style = $s_D11
$s_E11 = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.12%  '
$ws.Range("E11").Style = $s_E11
$s_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.38'
$ws.Range("D12").Style = $s_D12
$s_E12 = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -14.18%  '
$ws.Range("E12").Style = $s_E12
$s_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("D13").Style = $s_D13
$s_E13 = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("E13").Style = $s_E13
$s_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.84'
$ws.Range("D14").Style = $s_D14
$s_E14 = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("E14").Style = $s_E14
$s_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.494.84'
$ws.Range("D15").Style = $s_D15
$s_E15 = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("E15").Style = $s_E15
$s_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.853'
$ws.Range("D16").Style = $s_D16
$s_E16 = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("E16").Style = $s_E16
$s_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.21'
$ws.Range("D17").Style = $s_D17
$s_E17 = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.71%  '
$ws.Range("E17").Style = $s_E17
$s_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.155.85'
$ws.Range("D18").Style = $s_D18
$s_E18 = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.62%  '
$ws.Range("E18").Style = $s_E18
$s_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '40.979.57'
$ws.Range("D19").Style = $s_D19
$s_E19 = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("E19").Style = $s_E19
$s_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0938'
$ws.Range("D20").Style = $s_D20
$s_E20 = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.25%  '
$ws.Range("E20").Style = $s_E20
$s_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.07'
$ws.Range("D21").Style = $s_D21
$s_E21 = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.55%  '
$ws.Range("E21").Style = $s_E21
$s_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.38'
$ws.Range("D22").Style = $s_D22
$s_E22 = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.27%  '
$ws.Range("E22").Style = $s_E22
$s_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.58'
$ws.Range("D23").Style = $s_D23
$s_E23 = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.41%  '
$ws.Range("E23").Style = $s_E23
$s_E24 = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -6.26%  '
$ws.Range("E24").Style = $s_E24
$s_E25 = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E25").Style = $s_E25
$s_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.29'
$ws.Range("D26").Style = $s_D26
$s_E26 = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +8.03%  '
$ws.Range("E26").Style = $s_E26
$s_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.69'
$ws.Range("D27").Style = $s_D27
$s_E27 = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.25%  '
$ws.Range("E27").Style = $s_E27
$s_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.42'
$ws.Range("D28").Style = $s_D28
$s_E28 = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.47%  '
$ws.Range("E28").Style = $s_E28
$s_E29 = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.39%  '
$ws.Range("E29").Style = $s_E29
$s_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.58'
$ws.Range("D30").Style = $s_D30
$s_E30 = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("E30").Style = $s_E30
$s_E31 = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.31%  '
$ws.Range("E31").Style = $s_E31
$s_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.19'
$ws.Range("D32").Style = $s_D32
$s_E32 = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.41%  '
$ws.Range("E32").Style = $s_E32
$s_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.122'
$ws.Range("D33").Style = $s_D33
$s_E33 = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.22%  '
$ws.Range("E33").Style = $s_E33
$s_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.63'
$ws.Range("D34").Style = $s_D34
$s_E34 = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.43%  '
$ws.Range("E34").Style = $s_E34
$s_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0747'
$ws.Range("D35").Style = $s_D35
$s_E35 = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.51%  '
$ws.Range("E35").Style = $s_E35
$s_E36 = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.06%  '
$ws.Range("E36").Style = $s_E36
$s_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.55'
$ws.Range("D37").Style = $s_D37
$s_E37 = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.14%  '
$ws.Range("E37").Style = $s_E37
$s_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.00'
$ws.Range("D38").Style = $s_D38
$s_E38 = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("E38").Style = $s_E38
$s_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.29'
$ws.Range("D39").Style = $s_D39
$s_E39 = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.00%  '
$ws.Range("E39").Style = $s_E39
$s_E40 = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.11%  '
$ws.Range("E40").Style = $s_E40
$s_E41 = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.80%  '
$ws.Range("E41").Style = $s_E41
$s_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.46'
$ws.Range("D42").Style = $s_D42
$s_E42 = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.31%  '
$ws.Range("E42").Style = $s_E42
$s_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.88'
$ws.Range("D43").Style = $s_D43
$s_E43 = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("E43").Style = $s_E43
$s_B44 = $ws.Range("B44").Style
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("B44").Style = $s_B44
$s_C44 = $ws.Range("C44").Style
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("C44").Style = $s_C44
$s_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.37'
$ws.Range("D44").Style = $s_D44
$s_E44 = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -11.52%  '
$ws.Range("E44").Style = $s_E44
$s_B45 = $ws.Range("B45").Style
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Celestia'
$ws.Range("B45").Style = $s_B45
$s_C45 = $ws.Range("C45").Style
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("C45").Style = $s_C45
$s_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.21'
$ws.Range("D45").Style = $s_D45
$s_E45 = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -5.07%  '
$ws.Range("E45").Style = $s_E45
$s_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.191'
$ws.Range("D46").Style = $s_D46
$s_E46 = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.22%  '
$ws.Range("E46").Style = $s_E46
$s_E47 = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.69%  '
$ws.Range("E47").Style = $s_E47
$s_B48 = $ws.Range("B48").Style
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'BinanceUSD'
$ws.Range("B48").Style = $s_B48
$s_C48 = $ws.Range("C48").Style
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("C48").Style = $s_C48
$s_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").Style = $s_D48
$s_E48 = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("E48").Style = $s_E48
$s_B49 = $ws.Range("B49").Style
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Cronos'
$ws.Range("B49").Style = $s_B49
$s_C49 = $ws.Range("C49").Style
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C49").Style = $s_C49
$s_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0988'
$ws.Range("D49").Style = $s_D49
$s_E49 = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("E49").Style = $s_E49
$s_E50 = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("E50").Style = $s_E50
$s_E51 = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.87%  '
$ws.Range("E51").Style = $s_E51
